$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "usuario" value for the test data row from "pruebasregistro49" to "pruebauser01"
$ws.Range("D2").Value = "pruebauser01"
